$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.148.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.60%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.970.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.22%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'596.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.25%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  +2.65%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.12%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'2.969.64"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.21%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.510"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.04%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'7.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +5.52%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  +8.18%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.448"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.51%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.0000243"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +7.45%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'33.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.37%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("E15").Value = "'  -0.72%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'3.461.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.10%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'63.027.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.65%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'  +1.02%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'2.963.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.94%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'446.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.98%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'  +0.51%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.675"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.73%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'7.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.13%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'11.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +4.47%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'81.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.10%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'2.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.33%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'11.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.55%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  +0.12%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'2.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.50%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'0.0000107"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +21.46%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'7.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +4.95%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  +0.84%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'26.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.48%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  -0.26%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.16%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  +8.14%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.71%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  +0.43%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  +2.93%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'49.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.24%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  -0.26%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  -4.65%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.288"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.58%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'41.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.73%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'2.716.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.10%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'371.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.20%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  -1.50%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'135.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.96%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  +0.03%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'23.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.28%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.106"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.02%  "
$ws.Range("E51").Style = "Normal"

